# INS final commit remaining unstaged changes
# Updates the TabQuery text for the Programs / Projects / Grants / Publications
# tabs on Sheet1 to match the latest SQL used by the automation harness.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ProgramsTab (row 2, column B) ---------------------------------------
# "Focus Area" column renamed to "Special Topic" and the "Data Location
# Details" CASE now falls back to prg.program_acronym instead of prg.website.
$programsQuery = @"
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Special Topic",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.program_acronym     
        ELSE prg.data_link
    END AS "Data Location Details" 
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Liver Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
"@

# --- ProjectsTab (row 3, column B) ----------------------------------------
# "org_name" renamed to "project_org_name".
$projectsQuery = @"
SELECT DISTINCT
    prj.project_id AS "Project ID", 
    prj.project_title AS "Project Title",
    prj.project_org_name AS "Organization",
    prj.project_start_date AS "Project Start Date",
    prj.project_end_date AS "Project End Date"
FROM 
    df_project prj
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
     prg.cancer_type LIKE '%Liver Cancer%'
ORDER BY 
    lower(prj.project_id) ASC
LIMIT 100;
"@

# --- GrantsTab (row 4, column B) -------------------------------------------
# "gnt.project_end_date" renamed to "gnt.grant_end_date"; extra space added
# before the LIKE keyword in the WHERE clause.
$grantsQuery = @"
SELECT DISTINCT
    gnt.grant_id AS "Grant ID", 
    prj.project_id AS "Project",
    gnt.grant_title AS "Grant Title",
    gnt.principal_investigators AS "Principal Investigators",
    gnt.program_officers AS "Program Officers",
    gnt.fiscal_year AS "Fiscal Year",
    gnt.grant_end_date AS "Project End Date"
FROM 
    df_grant gnt
LEFT JOIN 
    df_project prj ON gnt."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type  LIKE '%Liver Cancer%'
ORDER BY 
    lower(gnt.grant_id) ASC
LIMIT 100;
"@

# --- PublicationsTab (row 5, column B) -------------------------------------
# "pub.title" renamed to "pub.publication_title"; new CASE branch added for
# relative_citation_ratio = 1.0; extra space added before LIKE.
$publicationsQuery = @"
SELECT DISTINCT
    pub.pmid AS "PubMed ID", 
    pub.publication_title AS "Title",
    pub.authors AS "Authors",
    pub.publication_date AS "Publication Date",
    pub.cited_by AS "Cited By",
    CASE 
    WHEN pub.relative_citation_ratio = 0 THEN '0'
    WHEN pub.relative_citation_ratio = 7.0 THEN '7'
    WHEN pub.relative_citation_ratio = 2.0 THEN '2'
  WHEN pub.relative_citation_ratio = 1.0 THEN '1'
    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) 
    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)
END AS "Relative Citation Ratio"
FROM 
    df_publication pub
LEFT JOIN 
    df_project prj ON pub."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
WHERE 
     prg.cancer_type  LIKE '%Liver Cancer%'
ORDER BY 
    lower(pub.pmid) ASC
LIMIT 100;
"@

$ws.Range("B2").Value = $programsQuery
$ws.Range("B3").Value = $projectsQuery
$ws.Range("B4").Value = $grantsQuery
$ws.Range("B5").Value = $publicationsQuery

# --- Normalize formatting on the TabQuery column -----------------------
# The refreshed workbook collapses the per-row duplicate "wrap text" fonts
# down to a single shared style; reapply a uniform font so all the query
# cells (and the StatQuery cell in C2) end up referencing the same style.
$queryRange = $ws.Range("B2:B5")
$queryRange.Font.Size = 12
$queryRange.WrapText = $true

$ws.Range("C2").Font.Size = 12
$ws.Range("C2").WrapText = $true

# --- Selection / view state -------------------------------------------------
$ws.Range("C2").Select()
